# Generate Report for Handoff
# - Bump the "Ready for handoff" rows' Latest-Handoff timestamp to the
#   newly generated value (Overview!G, zh-cn!H, de-de!H for rows 7,8,9,11,12,13).
# - Mark those same rows' Priority column ("ht") on the zh-cn / de-de sheets,
#   since the handoff file was just (re)generated for them.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" column G ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-18 04:21:48"
}

# --- zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-18 04:21:43"
}

# --- de-de sheet: Priority (E) + Latest Handoff Datetime (H) ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-18 04:21:48"
}
